$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '59.499.60'
$ws.Range('E2').Value = '  +2.80%  '
$ws.Range('D3').Value = '2.993.35'
$ws.Range('E3').Value = '  +2.00%  '
$ws.Range('E4').Value = '  -0.02%  '
$ws.Range('D5').Value = '563.67'
$ws.Range('E5').Value = '  +2.46%  '
$ws.Range('D6').Value = '138.95'
$ws.Range('E6').Value = '  +4.96%  '
$ws.Range('E7').Value = '  -0.06%  '
$ws.Range('D8').Value = '0.519'
$ws.Range('E8').Value = '  +1.24%  '
$ws.Range('D9').Value = '2.983.92'
$ws.Range('E9').Value = '  +1.91%  '
$ws.Range('E11').Value = '  +11.26%  '
$ws.Range('E12').Value = '  +1.94%  '
$ws.Range('E13').Value = '  +4.22%  '
$ws.Range('D14').Value = '33.80'
$ws.Range('E14').Value = '  +2.80%  '
$ws.Range('E15').Value = '  -0.32%  '
$ws.Range('D16').Value = '3.490.40'
$ws.Range('E16').Value = '  +2.07%  '
$ws.Range('E17').Value = '  +4.27%  '
$ws.Range('D18').Value = '2.992.21'
$ws.Range('E18').Value = '  +2.12%  '
$ws.Range('D19').Value = '59.490.43'
$ws.Range('E19').Value = '  +2.84%  '
$ws.Range('D20').Value = '435.11'
$ws.Range('E20').Value = '  +4.47%  '
$ws.Range('D21').Value = '13.59'
$ws.Range('E21').Value = '  +2.44%  '
$ws.Range('E22').Value = '  +3.57%  '
$ws.Range('D23').Value = '13.40'
$ws.Range('E23').Value = '  -0.38%  '
$ws.Range('D24').Value = '7.03'
$ws.Range('E24').Value = '  +0.72%  '
$ws.Range('D25').Value = '79.99'
$ws.Range('E25').Value = '  +0.68%  '
$ws.Range('E26').Value = '  -0.11%  '
$ws.Range('B27').Value = 'ImmutableX'
$ws.Range('C27').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D27').Value = '2.21'
$ws.Range('E27').Value = '  +9.64%  '
$ws.Range('B28').Value = 'FirstDigitalUSD'
$ws.Range('C28').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('D28').Value = '1.00'
$ws.Range('E28').Value = '  +0.02%  '
$ws.Range('E29').Value = '  +2.85%  '
$ws.Range('E30').Value = '  +4.84%  '
$ws.Range('D31').Value = '0.108'
$ws.Range('E31').Value = '  +10.25%  '
$ws.Range('D32').Value = '6.26'
$ws.Range('E32').Value = '  +4.73%  '
$ws.Range('D33').Value = '25.77'
$ws.Range('E33').Value = '  +2.05%  '
$ws.Range('D34').Value = '0.0₃0780'
$ws.Range('E34').Value = '  +12.10%  '
$ws.Range('D35').Value = '0.993'
$ws.Range('E35').Value = '  +5.90%  '
$ws.Range('E36').Value = '  +3.96%  '
$ws.Range('D37').Value = '2.09'
$ws.Range('E37').Value = '  +0.41%  '
$ws.Range('D38').Value = '48.83'
$ws.Range('E38').Value = '  +1.07%  '
$ws.Range('D39').Value = '8.66'
$ws.Range('E39').Value = '  -0.62%  '
$ws.Range('D40').Value = '2.77'
$ws.Range('E40').Value = '  +6.55%  '
$ws.Range('D41').Value = '401.70'
$ws.Range('E41').Value = '  +6.95%  '
$ws.Range('E42').Value = '  +2.80%  '
$ws.Range('D43').Value = '2.759.12'
$ws.Range('E43').Value = '  +2.39%  '
$ws.Range('D44').Value = '0.106'
$ws.Range('E44').Value = '  -1.86%  '
$ws.Range('E45').Value = '  +6.08%  '
$ws.Range('E46').Value = '  -0.04%  '
$ws.Range('D47').Value = '34.69'
$ws.Range('E47').Value = '  +19.33%  '
$ws.Range('D48').Value = '122.90'
$ws.Range('E48').Value = '  -0.59%  '
$ws.Range('E49').Value = '  +1.67%  '
$ws.Range('E50').Value = '  +2.81%  '
$ws.Range('D51').Value = '23.47'
$ws.Range('E51').Value = '  +2.22%  '
